# Linux Lab 2 Users and Groups -- "changes from Fall 23"
#
# The Reading section used to have four separate paragraphs:
#   1. "Read (or listen to) CyberAces Module 1 Linux, Session 4, ..." (+ its
#      "Scroll down to Users and Groups" sans.org hyperlink)
#   2. "Or download the PDF" (+ its assets.contentstack.io hyperlink)
#   3. "Read about su and sudo in "The Linux Command Line" ..." (+ its
#      linuxcommand.org hyperlink)
#   4. "Download a PDF of the book." (+ its sourceforge.net hyperlink)
#
# The edit removes paragraphs 1 and 2 entirely (the CyberAces reading is no
# longer assigned), so the "Read about su and sudo ..." paragraph becomes the
# first paragraph of the Reading section, immediately followed by the
# "Download a PDF of the book." paragraph.

$d = $word.ActiveDocument

# Find the two paragraphs that are being removed by matching on their
# distinctive leading text (robust to any paragraph-index drift).
$firstIdx = -1
$lastIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($firstIdx -eq -1 -and $ptext.StartsWith("Read (or listen to) CyberAces")) {
        $firstIdx = $i
    }
    if ($ptext.StartsWith("Or download the PDF")) {
        $lastIdx = $i
    }
}

if ($firstIdx -ne -1 -and $lastIdx -ne -1 -and $lastIdx -ge $firstIdx) {
    $startPara = $d.Paragraphs.Item($firstIdx)
    $endPara = $d.Paragraphs.Item($lastIdx)
    $killRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $killRange.Delete()
    Write-Output "Removed paragraphs $firstIdx..$lastIdx (CyberAces reading + PDF download)."
} else {
    Write-Output "WARNING: could not locate the CyberAces paragraphs (first=$firstIdx last=$lastIdx); no change made."
}
